# Relatório Atualizado: Sun Feb 15 06:35:16 UTC 2026
# Append the newest Performance snapshot row (row 4) and extend both
# charts' series ranges so they pick up the new data point.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Performance")

# --- Append the new snapshot row ------------------------------------------
$newRow = 4
$ws.Cells.Item($newRow, 1).Value = "15/02/2026 06:35:16"   # Data
$ws.Cells.Item($newRow, 2).Value = 60.08                    # Banca Total ($)
$ws.Cells.Item($newRow, 3).Value = 2.4                      # Investido ($)
$ws.Cells.Item($newRow, 4).Value = 57.67                    # Livre ($)
$ws.Cells.Item($newRow, 5).Value = 0.08                     # PnL Hoje ($)
$ws.Cells.Item($newRow, 6).Value = 1                        # Trades Hoje
$ws.Cells.Item($newRow, 7).Value = "GRID"                   # Modo

# --- Extend the charts so the new row is plotted too ----------------------
$co = $ws.ChartObjects()
for ($i = 1; $i -le $co.Count; $i++) {
    $chart = $co.Item($i).Chart
    $series = $chart.SeriesCollection()
    for ($s = 1; $s -le $series.Count; $s++) {
        $ser = $series.Item($s)
        if ($ser.Formula -match '^=SERIES\((.*)\)$') {
            $parts = $matches[1]

            # Split on commas that are not inside double quotes.
            $args = @()
            $depth = 0
            $inQuotes = $false
            $cur = ""
            foreach ($ch in $parts.ToCharArray()) {
                if ($ch -eq '"') { $inQuotes = -not $inQuotes }
                if ($ch -eq ',' -and -not $inQuotes) {
                    $args += $cur
                    $cur = ""
                } else {
                    $cur += $ch
                }
            }
            $args += $cur

            if ($args.Count -ge 3) {
                $nameArg = $args[0]
                $catArg = $args[1]
                $valArg = $args[2]
                $orderArg = ""
                if ($args.Count -ge 4) { $orderArg = "," + $args[3] }

                if ($catArg -match '^(.*!\$[A-Z]+\$)(\d+):(\$[A-Z]+\$)(\d+)$') {
                    $catArg = $matches[1] + $matches[2] + ":" + $matches[3] + $newRow
                }
                if ($valArg -match '^(.*!\$[A-Z]+\$)(\d+):(\$[A-Z]+\$)(\d+)$') {
                    $valArg = $matches[1] + $matches[2] + ":" + $matches[3] + $newRow
                }

                $ser.Formula = "=SERIES(" + $nameArg + "," + $catArg + "," + $valArg + $orderArg + ")"
            }
        }
    }
}
